# feat: Spawn entity request
# Adds "Bird" and "Snake" entity prefab entries to the "Prefabs View" sheet,
# mirroring the existing "Monkey" entry (Architecture ID / Prefab resource path).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prefabs View")

# Copy the formatting of the existing data row (row 2) down into the two
# new rows so the new cells keep the same style as the rest of the table.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C4").PasteSpecial(-4122)

# Row 3: Bird
$ws.Range("A3").Value = "Bird view"
$ws.Range("B3").Value = "Bird"
$ws.Range("C3").Value = "Prefabs/Entities/LivingEntities/Animals/Bird"

# Row 4: Snake
$ws.Range("A4").Value = "Snake view"
$ws.Range("B4").Value = "Snake"
$ws.Range("C4").Value = "Prefabs/Entities/LivingEntities/Animals/Snake"
